# Update the "想去人数" (want-to-go count) values in column F
# for both the "展览" and "全部类型" worksheets.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Mapping of row -> new value for column F
$updates = @{
    2  = 1507
    3  = 27
    4  = 971
    5  = 64
    6  = 2359
    8  = 1442
    10 = 157
    11 = 50
    12 = 393
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
